$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A9").Value = 4
$ws.Range("A14").Value = 5
$ws.Range("A16").Value = 6
$ws.Range("A17").Value = 6
$ws.Range("A18").Value = 7
$ws.Range("A22").Value = 8
$ws.Range("A23").Value = 8
$ws.Range("A24").Value = 9
$ws.Range("A25").Value = 10
$ws.Range("A26").Value = 10
$ws.Range("A27").Value = 10
$ws.Range("A28").Value = 11
$ws.Range("A29").Value = 11
$ws.Range("A30").Value = 11
$ws.Range("A31").Value = 11
$ws.Range("A32").Value = 12
$ws.Range("A33").Value = 12
$ws.Range("A34").Value = 12
$ws.Range("A35").Value = 12
$ws.Range("A36").Value = 12
$ws.Range("A37").Value = 13
$ws.Range("A38").Value = 13
$ws.Range("A39").Value = 13
$ws.Range("A40").Value = 13
$ws.Range("A41").Value = 14
$ws.Range("A42").Value = 14
$ws.Range("A43").Value = 14
$ws.Range("A44").Value = 15
$ws.Range("A45").Value = 15
$ws.Range("A46").Value = 16
$ws.Range("A47").Value = 16
$ws.Range("A48").Value = 16
$ws.Range("A49").Value = 17
$ws.Range("A50").Value = 17
$ws.Range("A54").Value = 18
$ws.Range("A55").Value = 18
$ws.Range("A56").Value = 18
$ws.Range("A57").Value = 19
$ws.Range("A58").Value = 20
$ws.Range("A59").Value = 20
$ws.Range("A60").Value = 20
$ws.Range("A61").Value = 20
$ws.Range("A62").Value = 20
$ws.Range("A63").Value = 20
$ws.Range("A64").Value = 21
$ws.Range("A65").Value = 21
$ws.Range("A66").Value = 21
$ws.Range("A67").Value = 21
$ws.Range("A68").Value = 21
$ws.Range("A69").Value = 21
$ws.Range("A70").Value = 21
$ws.Range("A71").Value = 22
$ws.Range("A72").Value = 22
$ws.Range("A73").Value = 22
$ws.Range("A74").Value = 23
$ws.Range("A75").Value = 23
$ws.Range("A76").Value = 24
$ws.Range("A77").Value = 24
$ws.Range("A78").Value = 25
$ws.Range("A79").Value = 25
$ws.Range("A80").Value = 25
$ws.Range("A81").Value = 25
$ws.Range("A84").Value = 26
$ws.Range("A85").Value = 26
$ws.Range("A86").Value = 27
$ws.Range("A87").Value = 27
$ws.Range("A88").Value = 27
$ws.Range("A89").Value = 28
$ws.Range("A90").Value = 28
$ws.Range("A91").Value = 28
$ws.Range("A93").Value = 29
$ws.Range("A94").Value = 29
$ws.Range("A97").Value = 29